# PnL_by_Day.xlsx — refresh of the "Query1" Power Query table (Sheet1!A1:C).
# A new day's row (BTCUSDT, 2022-02-21, -146.31) was pulled in at the top of
# the query result and the whole table shifted/re-rendered beneath it, so the
# full data block (rows 2-36) is rewritten to match the refreshed values, the
# table is resized from A1:C35 to A1:C36, and the ExternalData_1 defined name
# is updated to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- rewrite the data rows (Pair / date-serial / P&L) -----------------
$ws.Range("A2").Value = "BTCUSDT"
$ws.Range("B2").Value = 44613
$ws.Range("C2").Value = -146.31

$ws.Range("A3").Value = "BTCUSDT"
$ws.Range("B3").Value = 44612
$ws.Range("C3").Value = -102.28

$ws.Range("A4").Value = "BTCUSDT"
$ws.Range("B4").Value = 44611
$ws.Range("C4").Value = 41.31

$ws.Range("A5").Value = "BTCUSDT"
$ws.Range("B5").Value = 44610
$ws.Range("C5").Value = 337.4

$ws.Range("A6").Value = "BTCUSDT"
$ws.Range("B6").Value = 44609
$ws.Range("C6").Value = -166.98

$ws.Range("A7").Value = "BTCUSDT"
$ws.Range("B7").Value = 44608
$ws.Range("C7").Value = -158.44999999999999

$ws.Range("A8").Value = "BTCUSDT"
$ws.Range("B8").Value = 44606
$ws.Range("C8").Value = -103.18

$ws.Range("A9").Value = "BTCUSDT"
$ws.Range("B9").Value = 44605
$ws.Range("C9").Value = -11.55

$ws.Range("A10").Value = "BTCUSDT"
$ws.Range("B10").Value = 44604
$ws.Range("C10").Value = 260.14999999999998

$ws.Range("A11").Value = "BTCUSDT"
$ws.Range("B11").Value = 44603
$ws.Range("C11").Value = 252.52

$ws.Range("A12").Value = "BTCUSDT"
$ws.Range("B12").Value = 44602
$ws.Range("C12").Value = -244.53

$ws.Range("A13").Value = "BTCUSDT"
$ws.Range("B13").Value = 44601
$ws.Range("C13").Value = 5.25

$ws.Range("A14").Value = "ETHUSDT"
$ws.Range("B14").Value = 44587
$ws.Range("C14").Value = -290.8

$ws.Range("A15").Value = "ETHUSDT"
$ws.Range("B15").Value = 44586
$ws.Range("C15").Value = -831.64

$ws.Range("A16").Value = "BTCUSDT"
$ws.Range("B16").Value = 44585
$ws.Range("C16").Value = -45.64

$ws.Range("A17").Value = "BTCUSDT"
$ws.Range("B17").Value = 44585
$ws.Range("C17").Value = 14.26

$ws.Range("A18").Value = "ETHUSDT"
$ws.Range("B18").Value = 44584
$ws.Range("C18").Value = -387

$ws.Range("A19").Value = "ETHUSDT"
$ws.Range("B19").Value = 44583
$ws.Range("C19").Value = -874.83

$ws.Range("A20").Value = "ETHUSDT"
$ws.Range("B20").Value = 44582
$ws.Range("C20").Value = -1131.73

$ws.Range("A21").Value = "ETHUSDT"
$ws.Range("B21").Value = 44581
$ws.Range("C21").Value = 24.08

$ws.Range("A22").Value = "BTCUSDT"
$ws.Range("B22").Value = 44580
$ws.Range("C22").Value = 2.62

$ws.Range("A23").Value = "ETHUSDT"
$ws.Range("B23").Value = 44580
$ws.Range("C23").Value = -212.56

$ws.Range("A24").Value = "BTCUSDT"
$ws.Range("B24").Value = 44579
$ws.Range("C24").Value = -1832.65

$ws.Range("A25").Value = "BTCUSDT"
$ws.Range("B25").Value = 44578
$ws.Range("C25").Value = -104.96

$ws.Range("A26").Value = "ETHUSDT"
$ws.Range("B26").Value = 44577
$ws.Range("C26").Value = -30.73

$ws.Range("A27").Value = "BTCUSDT"
$ws.Range("B27").Value = 44577
$ws.Range("C27").Value = -103.03

$ws.Range("A28").Value = "BTCUSDT"
$ws.Range("B28").Value = 44576
$ws.Range("C28").Value = -85.29

$ws.Range("A29").Value = "BTCUSDT"
$ws.Range("B29").Value = 44575
$ws.Range("C29").Value = -18.149999999999999

$ws.Range("A30").Value = "BTCUSDT"
$ws.Range("B30").Value = 44574
$ws.Range("C30").Value = -43.46

$ws.Range("A31").Value = "BTCUSDT"
$ws.Range("B31").Value = 44573
$ws.Range("C31").Value = -1187.1300000000001

$ws.Range("A32").Value = "BTCUSDT"
$ws.Range("B32").Value = 44572
$ws.Range("C32").Value = -113.59

$ws.Range("A33").Value = "BTCUSDT"
$ws.Range("B33").Value = 44571
$ws.Range("C33").Value = -0.12

$ws.Range("A34").Value = "ETHUSDT"
$ws.Range("B34").Value = 44296
$ws.Range("C34").Value = -6.85

$ws.Range("A35").Value = "BTCUSDT"
$ws.Range("B35").Value = 44280
$ws.Range("C35").Value = -63.74

# New row 36 (was previously row 35's data, shifted down one by the refresh)
$ws.Range("A36").Value = "BTCUSDT"
$ws.Range("B36").Value = 44279
$ws.Range("C36").Value = 12.73

# --- grow the query table / autofilter from A1:C35 to A1:C36 ----------
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:C36"))

# --- keep the hidden ExternalData_1 defined name in sync ---------------
$wb.Names.Item("Sheet1!ExternalData_1").RefersTo = "=Sheet1!`$A`$1:`$C`$36"
